$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4628.1704
$ws.Range("I15").Value = 4628.1704
$ws.Range("K15").Value = 13884.5112
$ws.Range("M15").Value = -13715.5112

$ws.Range("H43").Value = 1754
$ws.Range("I43").Value = 1600
$ws.Range("J43").Value = 1761.7
$ws.Range("K43").Value = 1600
$ws.Range("L43").Value = 1761.7
$ws.Range("M43").Value = -1531
$ws.Range("N43").Value = -1899.7

$ws.Range("H46").Value = 433839
$ws.Range("I46").Value = 400508.5
$ws.Range("K46").Value = 1201525.5
$ws.Range("M46").Value = -1201406.5

$ws.Range("H60").Value = 433839
$ws.Range("I60").Value = 400508.5
$ws.Range("K60").Value = 1201525.5
$ws.Range("M60").Value = -1201041.5

$ws.Range("H137").Value = 2441321
$ws.Range("I137").Value = 4763729
$ws.Range("J137").Value = 2792.8
$ws.Range("K137").Value = 14291187
$ws.Range("L137").Value = 8378.400000000001
$ws.Range("M137").Value = -14288637
$ws.Range("N137").Value = -13478.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23544.338
$ws.Range("I32").Value = 20715.549
$ws.Range("J32").Value = 31896
$ws.Range("K32").Value = 20715.549
$ws.Range("L32").Value = 31896
$ws.Range("M32").Value = -20428.549
$ws.Range("N32").Value = -32470

$ws.Range("H63").Value = 1900
$ws.Range("I63").Value = 1900
$ws.Range("K63").Value = 1900
$ws.Range("M63").Value = -1214

$ws.Range("H66").Value = 1900
$ws.Range("I66").Value = 1900
$ws.Range("K66").Value = 9500
$ws.Range("M66").Value = -6068

$ws.Range("H74").Value = 10819101
$ws.Range("I74").Value = 14537637
$ws.Range("J74").Value = 128312.25
$ws.Range("K74").Value = 14537637
$ws.Range("L74").Value = 128312.25
$ws.Range("M74").Value = -14536763
$ws.Range("N74").Value = -130060.25

$ws.Range("H77").Value = 10819101
$ws.Range("I77").Value = 14537637
$ws.Range("J77").Value = 128312.25
$ws.Range("K77").Value = 72688185
$ws.Range("L77").Value = 641561.25
$ws.Range("M77").Value = -72683817
$ws.Range("N77").Value = -650297.25

$ws.Range("H80").Value = 34536.668
$ws.Range("J80").Value = 34444
$ws.Range("L80").Value = 34444
$ws.Range("N80").Value = -36440

$ws.Range("H83").Value = 34536.668
$ws.Range("J83").Value = 34444
$ws.Range("L83").Value = 103332
$ws.Range("N83").Value = -113316

$ws.Range("H122").Value = 3161.1667
$ws.Range("I122").Value = 2137.3333
$ws.Range("J122").Value = 4185
$ws.Range("K122").Value = 6411.999899999999
$ws.Range("L122").Value = 12555
$ws.Range("M122").Value = -3961.999899999999
$ws.Range("N122").Value = -17455

$ws.Range("H132").Value = 53139.293
$ws.Range("I132").Value = 37714.5
$ws.Range("J132").Value = 86361.92
$ws.Range("K132").Value = 113143.5
$ws.Range("L132").Value = 259085.76
$ws.Range("M132").Value = -110613.5
$ws.Range("N132").Value = -264145.76

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H123").Value = 24999.092
$ws.Range("J123").Value = 24999.092
$ws.Range("L123").Value = 24999.092
$ws.Range("N123").Value = -34799.092

$ws.Range("H134").Value = 2688.5
$ws.Range("I134").Value = 2487.8948
$ws.Range("K134").Value = 7463.6844
$ws.Range("M134").Value = -4928.6844

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2525.342
$ws.Range("I31").Value = 2133.0625
$ws.Range("J31").Value = 2810.6365
$ws.Range("K31").Value = 2133.0625
$ws.Range("L31").Value = 2810.6365
$ws.Range("M31").Value = -1838.0625
$ws.Range("N31").Value = -3400.6365

$ws.Range("H34").Value = 2525.342
$ws.Range("I34").Value = 2133.0625
$ws.Range("J34").Value = 2810.6365
$ws.Range("K34").Value = 2133.0625
$ws.Range("L34").Value = 2810.6365
$ws.Range("M34").Value = -1931.0625
$ws.Range("N34").Value = -3214.6365

$ws.Range("H94").Value = 3452.16
$ws.Range("I94").Value = 6559.7
$ws.Range("J94").Value = 1380.4667
$ws.Range("K94").Value = 6559.7
$ws.Range("L94").Value = 1380.4667
$ws.Range("M94").Value = -6108.7
$ws.Range("N94").Value = -2282.4667

$ws.Range("H122").Value = 2848.7778
$ws.Range("I122").Value = 2703.6667
$ws.Range("J122").Value = 2921.3333
$ws.Range("K122").Value = 8111.000100000001
$ws.Range("L122").Value = 8763.999899999999
$ws.Range("M122").Value = -5661.000100000001
$ws.Range("N122").Value = -13663.9999

$ws.Range("H132").Value = 18601.152
$ws.Range("I132").Value = 1320.1892
$ws.Range("J132").Value = 47664.59
$ws.Range("K132").Value = 3960.5676
$ws.Range("L132").Value = 142993.77
$ws.Range("M132").Value = -1430.5676
$ws.Range("N132").Value = -148053.77

$ws.Range("H134").Value = 50778
$ws.Range("I134").Value = 1135.6666
$ws.Range("J134").Value = 69393.875
$ws.Range("K134").Value = 3406.9998
$ws.Range("L134").Value = 208181.625
$ws.Range("M134").Value = -871.9998000000001
$ws.Range("N134").Value = -213251.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 560.2727
$ws.Range("I86").Value = 400
$ws.Range("J86").Value = 576.3
$ws.Range("K86").Value = 1200
$ws.Range("L86").Value = 1728.9
$ws.Range("N86").Value = -4100.9
$ws.Range("M86").Value = -14

$ws.Range("H89").Value = 560.2727
$ws.Range("I89").Value = 400
$ws.Range("J89").Value = 576.3
$ws.Range("K89").Value = 3600
$ws.Range("L89").Value = 5186.7
$ws.Range("N89").Value = -17042.7
$ws.Range("M89").Value = 2328

$ws.Range("H107").Value = 813.8929000000001
$ws.Range("J107").Value = 633.94116
$ws.Range("L107").Value = 1901.82348
$ws.Range("N107").Value = -5741.82348

$ws.Range("H113").Value = 574.75
$ws.Range("I113").Value = 465
$ws.Range("J113").Value = 640.6
$ws.Range("K113").Value = 1395
$ws.Range("L113").Value = 1921.8
$ws.Range("M113").Value = 775
$ws.Range("N113").Value = -6261.8

$ws.Range("H122").Value = 899.26086
$ws.Range("I122").Value = 354.66666
$ws.Range("J122").Value = 1249.3572
$ws.Range("K122").Value = 3191.99994
$ws.Range("L122").Value = 11244.2148
$ws.Range("M122").Value = -741.9999399999997
$ws.Range("N122").Value = -16144.2148

$ws.Range("H131").Value = 1170.8182
$ws.Range("I131").Value = 415.57144
$ws.Range("J131").Value = 2492.5
$ws.Range("K131").Value = 1246.71432
$ws.Range("L131").Value = 7477.5
$ws.Range("M131").Value = 3793.28568
$ws.Range("N131").Value = -17557.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 8667.333000000001
$ws.Range("J92").Value = 8667.333000000001
$ws.Range("L92").Value = 8667.333000000001
$ws.Range("N92").Value = -12411.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1080.5454
$ws.Range("I93").Value = 849.25
$ws.Range("K93").Value = 849.25
$ws.Range("M93").Value = 398.75

$ws.Range("H132").Value = 66600.81
$ws.Range("I132").Value = 4134.5835
$ws.Range("K132").Value = 12403.7505
$ws.Range("M132").Value = -9873.750499999998

$ws.Range("H136").Value = 87067
$ws.Range("I136").Value = 56142.26
$ws.Range("J136").Value = 204581
$ws.Range("K136").Value = 168426.78
$ws.Range("L136").Value = 613743
$ws.Range("M136").Value = -165876.78

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 85503.125
$ws.Range("I136").Value = 72700.21000000001
$ws.Range("J136").Value = 103427.2
$ws.Range("K136").Value = 218100.63
$ws.Range("L136").Value = 310281.6
$ws.Range("M136").Value = -215550.63
$ws.Range("N136").Value = -315381.6

Write-Output "All 33 row updates applied."
